$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Groups V2" sheet: scroll the view down (topLeftCell A166 -> A211).
#    Selection itself is left untouched (still A184:XFD184).
# ---------------------------------------------------------------------------
$wsV2 = $wb.Worksheets.Item("Groups V2")
$wsV2.Activate()
$excel.ActiveWindow.ScrollRow = 211
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------------
# 2. "Groups V3" sheet: the ortho_N group labels in column C get zero-padded
#    (ortho_0 -> ortho_00, ortho_1 -> ortho_01, ... ortho_9 -> ortho_09).
#    Written group-by-group, top to bottom, so the new shared strings land
#    in the workbook in ortho_00..ortho_09 order.
# ---------------------------------------------------------------------------
$wsV3 = $wb.Worksheets.Item("Groups V3")
$wsV3.Activate()

$wsV3.Range("C3:C6").Value = "ortho_00"
$wsV3.Range("C9:C22").Value = "ortho_01"
$wsV3.Range("C25:C34").Value = "ortho_02"
$wsV3.Range("C37:C47").Value = "ortho_03"
$wsV3.Range("C50:C60").Value = "ortho_04"
$wsV3.Range("C63:C75").Value = "ortho_05"
$wsV3.Range("C78:C90").Value = "ortho_06"
$wsV3.Range("C93:C109").Value = "ortho_07"
$wsV3.Range("C112:C133").Value = "ortho_08"
$wsV3.Range("C136").Value = "ortho_09"

# Move the live selection on this sheet from D9 to D23.
$wsV3.Range("D23").Select()

# ---------------------------------------------------------------------------
# 3. Add the new "Groups V3 Assessment" sheet right after "Groups V3" and
#    fill in its small header block.
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAssess = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$wsAssess.Name = "Groups V3 Assessment"

$wsAssess.Range("B3").Value = "RMSE / 50% / 75%"
$wsAssess.Range("B4").Value = "a"
$wsAssess.Range("C4").Value = "b"
$wsAssess.Range("D4").Value = "c"
$wsAssess.Range("E4").Value = "train"
$wsAssess.Range("F4").Value = "validation"
$wsAssess.Range("E3").Value = "Mean"
$wsAssess.Range("G3").Value = "Variance"
$wsAssess.Range("G4").Value = "train"
$wsAssess.Range("H4").Value = "validation"
$wsAssess.Range("E2").Value = "loss"
$wsAssess.Range("A4").Value = "group"

$wsAssess.Range("G4:H4").Select()
